$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Om Patel"

# Test Case 1 - __init__ valid input
$ws.Range("E7").Value = 'Valid client number, first_name, last_name, and email_address'
$ws.Range("F7").Value = 'Client(2904, "Om", "Patel", "ompatel@gmail.com")'
$ws.Range("G7").Value = 'Client object is successfully created with client_number=2904, first_name="Om", last_name="Patel", email_address="ompatel@gmail.com".'

# Test Case 2 - __init__ invalid client number
$ws.Range("E8").Value = 'Invalid client number (non-numeric).'
$ws.Range("F8").Value = 'Client("abc", "Om", "Patel", "ompatel@gmail.com")'
$ws.Range("G8").Value = 'ValueError: "Client number should be numeric."'

# Test Case 3 - __init__ blank first name
$ws.Range("E9").Value = 'Blank first_name (" ").'
$ws.Range("F9").Value = 'Client(1010, " ", "Patel", "ompatel@gmail.com")'
$ws.Range("G9").Value = 'ValueError: "First name cannot be blank."'

# Test Case 4 - __init__ blank last name
$ws.Range("E10").Value = 'Blank last_name (" ").'
$ws.Range("F10").Value = 'Client(1010, "Om", " ", "ompatel@gmail.com")'
$ws.Range("G10").Value = 'ValueError: "Last name cannot be blank."'

# Test Case 5 - __init__ invalid email
$ws.Range("E11").Value = 'Invalid email address.'
$ws.Range("F11").Value = 'Client(1010, "Om", "Patel", "invalid_email")'
$ws.Range("G11").Value = 'email_address is set to the default value "om.patel@pixell-river.com".'

# Test Case 6 - client_number getter
$ws.Range("E12").Value = 'Client object initialized.'
$ws.Range("F12").Value = 'client.client_number'
$ws.Range("G12").Value = 'Returns 2904.'

# Test Case 7 - first_name getter
$ws.Range("E13").Value = 'Client object initialized.'
$ws.Range("F13").Value = 'client.first_name'
$ws.Range("G13").Value = 'Returns "Om".'

# Test Case 8 - last_name getter
$ws.Range("E14").Value = 'Client object initialized.'
$ws.Range("F14").Value = 'client.last_name'
$ws.Range("G14").Value = 'Returns "Patel".'

# Test Case 9 - email_address getter
$ws.Range("E15").Value = 'Client object initialized.'
$ws.Range("F15").Value = 'client.email_address'
$ws.Range("G15").Value = 'Returns "ompatel@gmail.com".'

# Test Case 10 - __str__
$ws.Range("E16").Value = 'Client object initialized.'
$ws.Range("F16").Value = 'str(client)'
$ws.Range("G16").Value = 'Returns formatted string: "Name: Patel, Om\nClient Number: 2904\nEmail Address: ompatel@gmail.com\n".'

# Row 16 Preconditions/Inputs/Expected Result cells use the bold wrap style,
# matching rows 7-15, rather than the plain wrap style used by blank rows.
$ws.Range("E16:G16").Font.Bold = $true

# Rows grew taller because of the newly entered (wrapped) text.
$ws.Rows(7).RowHeight = 105
$ws.Rows(16).RowHeight = 106.8

# Restore the cursor/selection position recorded at save time.
$ws.Range("K17").Select()
